# Swap the values of columns C and D (codeforiati:group-code / codeforiati:group-name,
# and all the corresponding code/name pairs) for every row in the used range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)   # Column C
    $dCell = $ws.Cells.Item($r, 4)   # Column D

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
